$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q4
#    at the top of the data (row 2), pushing the rest down.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$summary.Application.CutCopyMode = $false

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 7
$summary.Cells.Item(2, 4).Value = 0.49

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q3"
$summary.Cells.Item(3, 3).Value = 8
$summary.Cells.Item(3, 4).Value = 0.42

$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = "2022-Q2"
$summary.Cells.Item(4, 3).Value = 7
$summary.Cells.Item(4, 4).Value = 0.06

$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(5, 2).Value = "2022-Q1"
$summary.Cells.Item(5, 3).Value = 5
$summary.Cells.Item(5, 4).Value = 0.11

$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(6, 2).Value = "2021-Q4"
$summary.Cells.Item(6, 3).Value = 5
$summary.Cells.Item(6, 4).Value = 0.21

$summary.Cells.Item(7, 1).Value = 5
$summary.Cells.Item(7, 2).Value = "2021-Q2"
$summary.Cells.Item(7, 3).Value = 1
$summary.Cells.Item(7, 4).Value = 0.04

# ------------------------------------------------------------------
# 2. Insert a brand new "2022-Q4" sheet right after "总计" containing
#    the quarterly fund holdings detail. Base it on a copy of the
#    "2022-Q3" sheet so formatting/column layout matches, then
#    overwrite the data with the 2022-Q4 numbers.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $summary)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template (old 2022-Q3) sheet has 8 data rows (rows 2-9); the new
# 2022-Q4 sheet only needs 7 data rows (rows 2-8), so drop the extra row.
$q4.Rows.Item(9).Delete()

function Set-TextCell($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Set-NumberCell($sheet, $row, $col, $number) {
    $sheet.Cells.Item($row, $col).Value = $number
}

# Row 2: 001521 国寿安保成长优选股票
Set-TextCell $q4 2 2 "001521"
Set-TextCell $q4 2 3 "国寿安保成长优选股票"
Set-TextCell $q4 2 4 "15.77"
Set-TextCell $q4 2 5 "88.89"
Set-TextCell $q4 2 6 "2.58"
Set-TextCell $q4 2 7 "0.4069"
Set-NumberCell $q4 2 8 8

# Row 3: 519929 长信电子信息行业量化灵活配置混合A
Set-TextCell $q4 3 2 "519929"
Set-TextCell $q4 3 3 "长信电子信息行业量化灵活配置混合A"
Set-TextCell $q4 3 4 "0.84"
Set-TextCell $q4 3 5 "90.40"
Set-TextCell $q4 3 6 "4.91"
Set-TextCell $q4 3 7 "0.0412"
Set-NumberCell $q4 3 8 7

# Row 4: 008082 国寿安保研究精选混合A
Set-TextCell $q4 4 2 "008082"
Set-TextCell $q4 4 3 "国寿安保研究精选混合A"
Set-TextCell $q4 4 4 "0.34"
Set-TextCell $q4 4 5 "84.56"
Set-TextCell $q4 4 6 "5.43"
Set-TextCell $q4 4 7 "0.0185"
Set-NumberCell $q4 4 8 1

# Row 5: 001226 中邮稳健添利灵活配置混合
Set-TextCell $q4 5 2 "001226"
Set-TextCell $q4 5 3 "中邮稳健添利灵活配置混合"
Set-TextCell $q4 5 4 "0.40"
Set-TextCell $q4 5 5 "93.24"
Set-TextCell $q4 5 6 "2.87"
Set-TextCell $q4 5 7 "0.0115"
Set-NumberCell $q4 5 8 6

# Row 6: 008083 国寿安保研究精选混合C
Set-TextCell $q4 6 2 "008083"
Set-TextCell $q4 6 3 "国寿安保研究精选混合C"
Set-TextCell $q4 6 4 "0.13"
Set-TextCell $q4 6 5 "84.56"
Set-TextCell $q4 6 6 "5.43"
Set-TextCell $q4 6 7 "0.0071"
Set-NumberCell $q4 6 8 1

# Row 7: 005536 渤海汇金量化成长混合
Set-TextCell $q4 7 2 "005536"
Set-TextCell $q4 7 3 "渤海汇金量化成长混合"
Set-TextCell $q4 7 4 "0.37"
Set-TextCell $q4 7 5 "86.00"
Set-TextCell $q4 7 6 "1.46"
Set-TextCell $q4 7 7 "0.0054"
Set-NumberCell $q4 7 8 3

# Row 8: 013153 长信电子信息行业量化灵活配置混合C (持有市值 column is a plain 0)
Set-TextCell $q4 8 2 "013153"
Set-TextCell $q4 8 3 "长信电子信息行业量化灵活配置混合C"
Set-TextCell $q4 8 4 "0.00"
Set-TextCell $q4 8 5 "90.40"
Set-TextCell $q4 8 6 "4.91"
Set-NumberCell $q4 8 7 0
Set-NumberCell $q4 8 8 7

# ------------------------------------------------------------------
# 3. Restore the originally-selected tab (the last sheet, "2021-Q2")
#    as the active sheet, since copying/inserting sheets above moved
#    the selection.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
